$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 296
$ws.Range('B296').Value = 6962998
$ws.Range('F296').Value = 'Bandirmaspor'
$ws.Range('G296').Value = 'Adanaspor'
$ws.Range('H296').Value = 2
$ws.Range('I296').Value = 3
$ws.Range('J296').Value = 'A'
$ws.Range('K296').Value = 1.4
$ws.Range('L296').Value = 4.5
$ws.Range('M296').Value = 7.5
$ws.Range('N296').Value = 1.4
$ws.Range('O296').Value = 4.75
$ws.Range('P296').Value = 7
$ws.Range('Q296').Value = -1.25
$ws.Range('R296').Value = 1.85
$ws.Range('S296').Value = 1.95
$ws.Range('T296').Value = 2.75
$ws.Range('U296').Value = 1.8
$ws.Range('V296').Value = 2
$ws.Range('W296').Value = -1
$ws.Range('Y296').Value = 6
$ws.Range('Z296').Value = -1
$ws.Range('AA296').Value = 0.95
$ws.Range('AB296').Value = 0.8

# Row 297
$ws.Range('B297').Value = 6962933
$ws.Range('F297').Value = 'Bodrum BLD Spor'
$ws.Range('G297').Value = 'Goztepe'
$ws.Range('H297').Value = 3
$ws.Range('I297').Value = 0
$ws.Range('J297').Value = 'H'
$ws.Range('K297').Value = 2.3
$ws.Range('L297').Value = 3.25
$ws.Range('M297').Value = 3.1
$ws.Range('N297').Value = 2.8
$ws.Range('O297').Value = 3.25
$ws.Range('P297').Value = 2.5
$ws.Range('Q297').Value = 0
$ws.Range('R297').Value = 1.975
$ws.Range('S297').Value = 1.825
$ws.Range('T297').Value = 2.25
$ws.Range('U297').Value = 1.9
$ws.Range('V297').Value = 1.9
$ws.Range('W297').Value = 1.8
$ws.Range('Y297').Value = -1
$ws.Range('Z297').Value = 0.9750000000000001
$ws.Range('AA297').Value = -1
$ws.Range('AB297').Value = 0.8999999999999999

# Row 335
$ws.Range('B335').Value = 6963045
$ws.Range('F335').Value = 'Erzurum BB'
$ws.Range('G335').Value = 'Sanliurfaspor'
$ws.Range('I335').Value = 0
$ws.Range('K335').Value = 1.95
$ws.Range('L335').Value = 3.5
$ws.Range('M335').Value = 3.75
$ws.Range('N335').Value = 1.95
$ws.Range('O335').Value = 3.5
$ws.Range('P335').Value = 3.75
$ws.Range('Q335').Value = -0.5
$ws.Range('R335').Value = 1.95
$ws.Range('S335').Value = 1.85
$ws.Range('T335').Value = 2.25
$ws.Range('U335').Value = 1.875
$ws.Range('V335').Value = 1.925
$ws.Range('W335').Value = 0.95
$ws.Range('Z335').Value = 0.95
$ws.Range('AB335').Value = -0.5
$ws.Range('AC335').Value = 0.4625

# Row 336
$ws.Range('B336').Value = 6963229
$ws.Range('F336').Value = 'Tuzlaspor'
$ws.Range('G336').Value = 'Manisa BBSK'
$ws.Range('I336').Value = 1
$ws.Range('K336').Value = 3.6
$ws.Range('L336').Value = 3.4
$ws.Range('M336').Value = 2.05
$ws.Range('N336').Value = 3.5
$ws.Range('O336').Value = 3.4
$ws.Range('P336').Value = 2.05
$ws.Range('Q336').Value = 0.25
$ws.Range('R336').Value = 2
$ws.Range('S336').Value = 1.8
$ws.Range('T336').Value = 2.5
$ws.Range('U336').Value = 2
$ws.Range('V336').Value = 1.8
$ws.Range('W336').Value = 2.5
$ws.Range('Z336').Value = 1
$ws.Range('AB336').Value = 1
$ws.Range('AC336').Value = -1

# Row 361
$ws.Range('B361').Value = 6963238
$ws.Range('F361').Value = 'Boluspor'
$ws.Range('G361').Value = 'Erzurum BB'
$ws.Range('K361').Value = 2.15
$ws.Range('L361').Value = 3.2
$ws.Range('M361').Value = 3.5
$ws.Range('N361').Value = 2.3
$ws.Range('O361').Value = 3.1
$ws.Range('P361').Value = 3.2
$ws.Range('Q361').Value = -0.25
$ws.Range('R361').Value = 1.925
$ws.Range('S361').Value = 1.875
$ws.Range('U361').Value = 2.025
$ws.Range('V361').Value = 1.775

# Row 362
$ws.Range('B362').Value = 6962946
$ws.Range('F362').Value = 'Keciorengucu'
$ws.Range('G362').Value = 'Kocaelispor'
$ws.Range('K362').Value = 3.8
$ws.Range('L362').Value = 3.4
$ws.Range('M362').Value = 1.95
$ws.Range('N362').Value = 4.2
$ws.Range('O362').Value = 3.5
$ws.Range('P362').Value = 1.85
$ws.Range('Q362').Value = 0.5
$ws.Range('R362').Value = 1.9
$ws.Range('S362').Value = 1.9
$ws.Range('U362').Value = 1.875
$ws.Range('V362').Value = 1.925
